$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.293.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.576.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.30%  "
$ws.Range("E9").Value = "  +4.08%  "
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.036.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.263.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.585.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "342.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("E20").Value = "  +3.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.695.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("E25").Value = "  +3.29%  "
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.990"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0828"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "468.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.93%  "
$ws.Range("E34").Value = "  +4.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "176.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.407"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.42%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "151.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("E43").Value = "  +2.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0551"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.616"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0982"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.78%  "
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.165"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.47%  "
